# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G) for rows 2-15, recalculated from the updated source data
$kValues = @{
    2  = 6
    3  = 7
    4  = 9
    5  = 1
    6  = 4
    7  = 3
    8  = 6
    9  = 4
    10 = 2
    11 = 4
    12 = 3
    13 = 6
    14 = 3
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
